$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly report (week ending 2022-07-04, serial 44746) is inserted at
# the top of this product's data block (rows 311/312), pushing the
# existing history down by two rows (old 311..328 -> new 313..330).
$ws.Rows("311:312").Insert()

# Duplicate the now-shifted rows 313:314 into the freshly inserted blank
# rows 311:312 so every "static" column (A,B,C,E,F,G,H,I,N,O,Q,R) is
# populated with the same market/category/quality metadata, matching the
# existing weekly pattern for this product.
$ws.Range("A313:R314").Copy()
$ws.Range("A311").PasteSpecial()

# Row 311 ("Primera" quality) - this week's reported figures
$ws.Range("D311").Value = 44746
$ws.Range("J311").Value = 160
$ws.Range("K311").Value = 14000
$ws.Range("L311").Value = 15000
$ws.Range("M311").Value = 14500
$ws.Range("P311").Value = 207

# Row 312 ("Segunda" quality) - this week's reported figures
$ws.Range("D312").Value = 44746
$ws.Range("J312").Value = 160
$ws.Range("K312").Value = 10000
$ws.Range("L312").Value = 11000
$ws.Range("M312").Value = 10500
$ws.Range("P312").Value = 105
